$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.34740256217261
$ws.Range("C2").Value = 10.34926009556185
$ws.Range("D2").Value = 9.964170495555722
$ws.Range("F2").Value = 30.15515580384336
$ws.Range("G2").Value = 29.67589848840167
$ws.Range("H2").Value = 14.55382193461381
$ws.Range("I2").Value = 22.77292271472182
$ws.Range("J2").Value = 10.53619613863462
$ws.Range("L2").Value = 11.73649658738849
$ws.Range("M2").Value = 15.61546650883949
$ws.Range("N2").Value = 18.16391667204133
$ws.Range("O2").Value = 22.27034709768229
$ws.Range("B3").Value = 13.92697751076037
$ws.Range("C3").Value = 10.20472931831081
$ws.Range("D3").Value = 9.974302547796395
$ws.Range("F3").Value = 30.21936931795404
$ws.Range("G3").Value = 29.70917269403619
$ws.Range("H3").Value = 14.594851974642
$ws.Range("I3").Value = 22.86734269547867
$ws.Range("J3").Value = 10.55588438663188
$ws.Range("L3").Value = 11.74004598782049
$ws.Range("M3").Value = 15.52957484803906
$ws.Range("N3").Value = 18.20351729125601
$ws.Range("O3").Value = 22.33051725549161
$ws.Range("B4").Value = 13.66353389365985
$ws.Range("C4").Value = 10.11432434482128
$ws.Range("D4").Value = 9.981694766624308
$ws.Range("F4").Value = 30.26537507066767
$ws.Range("G4").Value = 29.73852719587957
$ws.Range("H4").Value = 14.62225852458104
$ws.Range("I4").Value = 22.92932741198176
$ws.Range("J4").Value = 10.56860333210849
$ws.Range("L4").Value = 11.74350065000834
$ws.Range("M4").Value = 15.47829641035831
$ws.Range("N4").Value = 18.22948961852674
$ws.Range("O4").Value = 22.3719993008228
$ws.Range("B5").Value = 13.55500132588529
$ws.Range("C5").Value = 10.0770923071551
$ws.Range("D5").Value = 9.985002198086482
$ws.Range("F5").Value = 30.28577473776429
$ws.Range("G5").Value = 29.75272901947364
$ws.Range("H5").Value = 14.63398367033845
$ws.Range("I5").Value = 22.95559520973011
$ws.Range("J5").Value = 10.57394533808145
$ws.Range("L5").Value = 11.74523013571867
$ws.Range("M5").Value = 15.45778182351755
$ws.Range("N5").Value = 18.24049107711344
$ws.Range("O5").Value = 22.39004291853018
$ws.Range("B6").Value = 13.53691324663057
$ws.Range("C6").Value = 10.07088711809225
$ws.Range("D6").Value = 9.985569230069375
$ws.Range("F6").Value = 30.28926178812104
$ws.Range("G6").Value = 29.7552223363382
$ws.Range("H6").Value = 14.63596424946561
$ws.Range("I6").Value = 22.96001787208674
$ws.Range("J6").Value = 10.57484198627164
$ws.Range("L6").Value = 11.74553677220793
$ws.Range("M6").Value = 15.45439887469393
$ws.Range("N6").Value = 18.24234310092644
$ws.Range("O6").Value = 22.39310781192266
$ws.Range("B7").Value = 13.66207473781606
$ws.Range("C7").Value = 10.11382376897953
$ws.Range("D7").Value = 9.981738176464791
$ws.Range("F7").Value = 30.26564350225697
$ws.Range("G7").Value = 29.73870966521938
$ws.Range("H7").Value = 14.62241439967475
$ws.Range("I7").Value = 22.92967758511989
$ws.Range("J7").Value = 10.56867473216725
$ws.Range("L7").Value = 11.74352267070457
$ws.Range("M7").Value = 15.47801817839507
$ws.Range("N7").Value = 18.22963629648703
$ws.Range("O7").Value = 22.37223803240162
$ws.Range("B8").Value = 14.20363026225602
$ws.Range("C8").Value = 10.29978678604099
$ws.Range("D8").Value = 9.967421311566735
$ws.Range("F8").Value = 30.17593003268381
$ws.Range("G8").Value = 29.68551699536378
$ws.Range("H8").Value = 14.56750961803391
$ws.Range("I8").Value = 22.80464638487969
$ws.Range("J8").Value = 10.542854098994
$ws.Range("L8").Value = 11.73745634153055
$ws.Range("M8").Value = 15.58555719796905
$ws.Range("N8").Value = 18.17722743064953
$ws.Range("O8").Value = 22.290151123287
$ws.Range("B9").Value = 15.21728208675408
$ws.Range("C9").Value = 10.65012994927535
$ws.Range("D9").Value = 9.948612711769995
$ws.Range("F9").Value = 30.05228596706955
$ws.Range("G9").Value = 29.65215479256328
$ws.Range("H9").Value = 14.47740605472432
$ws.Range("I9").Value = 22.59127899534925
$ws.Range("J9").Value = 10.49720056964476
$ws.Range("L9").Value = 11.73563422938749
$ws.Range("M9").Value = 15.80735020876
$ws.Range("N9").Value = 18.08757056908049
$ws.Range("O9").Value = 22.16524758941806
$ws.Range("B10").Value = 15.92494749216199
$ws.Range("C10").Value = 10.89732488299516
$ws.Range("D10").Value = 9.940407927829662
$ws.Range("F10").Value = 29.99341765675967
$ws.Range("G10").Value = 29.67100493935721
$ws.Range("H10").Value = 14.42191154570722
$ws.Range("I10").Value = 22.45391112546772
$ws.Range("J10").Value = 10.46666659005507
$ws.Range("L10").Value = 11.74037412479661
$ws.Range("M10").Value = 15.97602941834582
$ws.Range("N10").Value = 18.02964984191868
$ws.Range("O10").Value = 22.09555406756461
$ws.Range("B11").Value = 16.23745627098432
$ws.Range("C11").Value = 11.0072718316443
$ws.Range("D11").Value = 9.937886071182088
$ws.Range("F11").Value = 29.97359264189896
$ws.Range("G11").Value = 29.68899003922769
$ws.Range("H11").Value = 14.39898917857355
$ws.Range("I11").Value = 22.39562971897274
$ws.Range("J11").Value = 10.45342287852901
$ws.Range("L11").Value = 11.74383568660323
$ws.Range("M11").Value = 16.05380633864379
$ws.Range("N11").Value = 18.00501692835259
$ws.Range("O11").Value = 22.06865671858129
$ws.Range("B12").Value = 16.35434170124182
$ws.Range("C12").Value = 11.04852260729425
$ws.Range("D12").Value = 9.937104404214287
$ws.Range("F12").Value = 29.9670857638036
$ws.Range("G12").Value = 29.69715069889009
$ws.Range("H12").Value = 14.39064296112507
$ws.Range("I12").Value = 22.3741653200752
$ws.Range("J12").Value = 10.44850031807887
$ws.Range("L12").Value = 11.74533282212803
$ws.Range("M12").Value = 16.08339096814111
$ws.Range("N12").Value = 17.99593506153023
$ws.Range("O12").Value = 22.05916352690858
$ws.Range("B13").Value = 16.32923454293206
$ws.Range("C13").Value = 11.03965594933916
$ws.Range("D13").Value = 9.937265055098397
$ws.Range("F13").Value = 29.96844264115176
$ws.Range("G13").Value = 29.6953331738787
$ws.Range("H13").Value = 14.39242561395809
$ws.Range("I13").Value = 22.37876112557945
$ws.Range("J13").Value = 10.44955636986092
$ws.Range("L13").Value = 11.74500212362133
$ws.Range("M13").Value = 16.07701379632181
$ws.Range("N13").Value = 17.99788006895232
$ws.Range("O13").Value = 22.06117725249712
$ws.Range("B14").Value = 16.24710220744642
$ws.Range("C14").Value = 11.01067338110569
$ws.Range("D14").Value = 9.937818294233683
$ws.Range("F14").Value = 29.97303726658953
$ws.Range("G14").Value = 29.6896343866379
$ws.Range("H14").Value = 14.3982958365492
$ws.Range("I14").Value = 22.39385169227948
$ws.Range("J14").Value = 10.45301604403433
$ws.Range("L14").Value = 11.74395513099575
$ws.Range("M14").Value = 16.05623773768697
$ws.Range("N14").Value = 18.0042648292592
$ws.Range("O14").Value = 22.06786182878557
$ws.Range("B15").Value = 16.19660148124893
$ws.Range("C15").Value = 10.99287005820921
$ws.Range("D15").Value = 9.938179715079295
$ws.Range("F15").Value = 29.97598189555136
$ws.Range("G15").Value = 29.6863194202196
$ws.Range("H15").Value = 14.40193501282792
$ws.Range("I15").Value = 22.40317396741845
$ws.Range("J15").Value = 10.45514723477656
$ws.Range("L15").Value = 11.74333803925365
$ws.Range("M15").Value = 16.04352848333607
$ws.Range("N15").Value = 18.00820770880451
$ws.Range("O15").Value = 22.07204650483803
$ws.Range("B16").Value = 15.90432600680784
$ws.Range("C16").Value = 10.89008716231917
$ws.Range("D16").Value = 9.940597031561795
$ws.Range("F16").Value = 29.99485325291141
$ws.Range("G16").Value = 29.67001869787167
$ws.Range("H16").Value = 14.42345631136159
$ws.Range("I16").Value = 22.45780466104888
$ws.Range("J16").Value = 10.46754506979249
$ws.Range("L16").Value = 11.74017405340783
$ws.Range("M16").Value = 15.97096599244756
$ws.Range("N16").Value = 18.03129412919547
$ws.Range("O16").Value = 22.09740868897428
$ws.Range("B17").Value = 15.72253834996252
$ws.Range("C17").Value = 10.8263747855167
$ws.Range("D17").Value = 9.942389492544754
$ws.Range("F17").Value = 30.00821184343208
$ws.Range("G17").Value = 29.6624271746265
$ws.Range("H17").Value = 14.43725377238543
$ws.Range("I17").Value = 22.49239685512745
$ws.Range("J17").Value = 10.4753160007526
$ws.Range("L17").Value = 11.73856641677723
$ws.Range("M17").Value = 15.92670593677658
$ws.Range("N17").Value = 18.04589584955726
$ws.Range("O17").Value = 22.11419946811568
$ws.Range("B18").Value = 15.6170983974007
$ws.Range("C18").Value = 10.78949484902305
$ws.Range("D18").Value = 9.94353445933975
$ws.Range("F18").Value = 30.01654994565294
$ws.Range("G18").Value = 29.6589467304149
$ws.Range("H18").Value = 14.44540829795978
$ws.Range("I18").Value = 22.5126894137107
$ws.Range("J18").Value = 10.47984649572256
$ws.Range("L18").Value = 11.7377647011269
$ws.Range("M18").Value = 15.90134806569415
$ws.Range("N18").Value = 18.05445588075055
$ws.Range("O18").Value = 22.1243094460384
$ws.Range("B19").Value = 15.58125027844051
$ws.Range("C19").Value = 10.77696844320139
$ws.Range("D19").Value = 9.943941726117792
$ws.Range("F19").Value = 30.01948549189051
$ws.Range("G19").Value = 29.65792055180151
$ws.Range("H19").Value = 14.4482068221454
$ws.Range("I19").Value = 22.51962813423516
$ws.Range("J19").Value = 10.48139090861549
$ws.Range("L19").Value = 11.73751441276432
$ws.Range("M19").Value = 15.89277993120371
$ws.Range("N19").Value = 18.05738191804994
$ws.Range("O19").Value = 22.12781016953266
$ws.Range("B20").Value = 15.74198185303579
$ws.Range("C20").Value = 10.83318149768466
$ws.Range("D20").Value = 9.942186889842858
$ws.Range("F20").Value = 30.00672204639957
$ws.Range("G20").Value = 29.66314362507135
$ws.Range("H20").Value = 14.4357623854727
$ws.Range("I20").Value = 22.4886734651239
$ws.Range("J20").Value = 10.4744824758983
$ws.Range("L20").Value = 11.7387248369793
$ws.Range("M20").Value = 15.93140734289974
$ws.Range("N20").Value = 18.0443247604735
$ws.Range("O20").Value = 22.11236523450648
$ws.Range("B21").Value = 16.27126670682776
$ws.Range("C21").Value = 11.01919686301296
$ws.Range("D21").Value = 9.937651097379817
$ws.Range("F21").Value = 29.9716605617516
$ws.Range("G21").Value = 29.6912716493385
$ws.Range("H21").Value = 14.39656254543377
$ws.Range("I21").Value = 22.38940279039034
$ws.Range("J21").Value = 10.4519973453605
$ws.Range("L21").Value = 11.74425761296085
$ws.Range("M21").Value = 16.06233672551969
$ws.Range("N21").Value = 18.00238279653262
$ws.Range("O21").Value = 22.06587961160181
$ws.Range("B22").Value = 16.60866506536073
$ws.Range("C22").Value = 11.1385216055567
$ws.Range("D22").Value = 9.935696463024666
$ws.Range("F22").Value = 29.95457694297149
$ws.Range("G22").Value = 29.71752182880822
$ws.Range("H22").Value = 14.37288990970786
$ws.Range("I22").Value = 22.32805322350939
$ws.Range("J22").Value = 10.43784125746388
$ws.Range("L22").Value = 11.74895890510306
$ws.Range("M22").Value = 16.1486687796367
$ws.Range("N22").Value = 17.97640538633952
$ws.Range("O22").Value = 22.0395339005382
$ws.Range("B23").Value = 16.42939963920196
$ws.Range("C23").Value = 11.0750489313768
$ws.Range("D23").Value = 9.936647559882088
$ws.Range("F23").Value = 29.96316124153017
$ws.Range("G23").Value = 29.70279319409575
$ws.Range("H23").Value = 14.38534630914337
$ws.Range("I23").Value = 22.36047354375973
$ws.Range("J23").Value = 10.44534741777288
$ws.Range("L23").Value = 11.74635090994485
$ws.Range("M23").Value = 16.10252796118606
$ws.Range("N23").Value = 17.99013899675169
$ws.Range("O23").Value = 22.05322556403517
$ws.Range("B24").Value = 15.73319432738353
$ws.Range("C24").Value = 10.83010496263417
$ws.Range("D24").Value = 9.942278129956017
$ws.Range("F24").Value = 30.00739353384395
$ws.Range("G24").Value = 29.66281696395949
$ws.Range("H24").Value = 14.43643594941744
$ws.Range("I24").Value = 22.49035554808801
$ws.Range("J24").Value = 10.47485911658039
$ws.Range("L24").Value = 11.73865283341876
$ws.Range("M24").Value = 15.92928156226227
$ws.Range("N24").Value = 18.04503453488014
$ws.Range("O24").Value = 22.11319306890888
$ws.Range("B25").Value = 14.94906624614994
$ws.Range("C25").Value = 10.55704167104207
$ws.Range("D25").Value = 9.952712355824813
$ws.Range("F25").Value = 30.08012574216012
$ws.Range("G25").Value = 29.65356779545498
$ws.Range("H25").Value = 14.49990133443883
$ws.Range("I25").Value = 22.64559470380076
$ws.Range("J25").Value = 10.5090208230985
$ws.Range("L25").Value = 11.74333803925365
$ws.Range("M25").Value = 15.74627519500859
$ws.Range("N25").Value = 18.00820770880451
$ws.Range("O25").Value = 22.07204650483803
